$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 999.6667
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 999.6667
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 999.6667
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1651.6667
$ws.Range("H98").Value = 1521.9375
$ws.Range("I98").Value = 1498.1072
$ws.Range("K98").Value = 1498.1072
$ws.Range("M98").Value = -0.1071999999999207
$ws.Range("H101").Value = 193.33333
$ws.Range("I101").Value = 193.33333
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 579.99999
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 1042.00001
$ws.Range("N101").ClearContents()
$ws.Range("H122").Value = 1521.9375
$ws.Range("I122").Value = 1498.1072
$ws.Range("K122").Value = 4494.321599999999
$ws.Range("M122").Value = -2044.321599999999
$ws.Range("H123").Value = 99999
$ws.Range("J123").Value = 99999
$ws.Range("L123").Value = 99999
$ws.Range("N123").Value = -109799
$ws.Range("H129").Value = 3219.6155
$ws.Range("I129").Value = 1181.5
$ws.Range("K129").Value = 3544.5
$ws.Range("M129").Value = 1455.5
$ws.Range("H131").Value = 1749482.8
$ws.Range("I131").Value = 986.4286
$ws.Range("K131").Value = 2959.2858
$ws.Range("M131").Value = 2080.7142
$ws.Range("H132").Value = 4349.625
$ws.Range("I132").Value = 2926.9092
$ws.Range("K132").Value = 8780.7276
$ws.Range("M132").Value = -6250.7276
$ws.Range("H137").Value = 542847.0600000001
$ws.Range("J137").Value = 1084501.8
$ws.Range("L137").Value = 3253505.4
$ws.Range("N137").Value = -3258605.4
$ws.Range("H138").Value = 2497.1128
$ws.Range("I138").Value = 954.4167
$ws.Range("J138").Value = 4083.8857
$ws.Range("K138").Value = 2863.2501
$ws.Range("L138").Value = 12251.6571
$ws.Range("M138").Value = 2276.7499
$ws.Range("N138").Value = -22531.6571

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 299
$ws.Range("I5").Value = 315.27274
$ws.Range("K5").Value = 315.27274
$ws.Range("M5").Value = -203.27274
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H44").Value = 89999
$ws.Range("J44").Value = 89999
$ws.Range("L44").Value = 89999
$ws.Range("N44").Value = -90975
$ws.Range("H60").Value = 70570.07000000001
$ws.Range("I60").Value = 74079.30499999999
$ws.Range("J60").Value = 24950
$ws.Range("K60").Value = 74079.30499999999
$ws.Range("L60").Value = 24950
$ws.Range("M60").Value = -73346.30499999999
$ws.Range("N60").Value = -26416
$ws.Range("H61").Value = 4672578.5
$ws.Range("I61").Value = 5411459
$ws.Range("J61").Value = 1255254
$ws.Range("K61").Value = 5411459
$ws.Range("L61").Value = 1255254
$ws.Range("M61").Value = -5411247
$ws.Range("N61").Value = -1255678
$ws.Range("H63").Value = 2933.875
$ws.Range("I63").Value = 2924.4285
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 2924.4285
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -2238.4285
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2933.875
$ws.Range("I66").Value = 2924.4285
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 14622.1425
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -11190.1425
$ws.Range("N66").Value = -21864
$ws.Range("H110").Value = 5473.3687
$ws.Range("I110").Value = 5882.1177
$ws.Range("K110").Value = 5882.1177
$ws.Range("M110").Value = -3837.1177
$ws.Range("H132").Value = 1855198
$ws.Range("I132").Value = 3259.7551
$ws.Range("K132").Value = 9779.265299999999
$ws.Range("M132").Value = -7249.265299999999
$ws.Range("H136").Value = 4672578.5
$ws.Range("I136").Value = 5411459
$ws.Range("J136").Value = 1255254
$ws.Range("K136").Value = 16234377
$ws.Range("L136").Value = 3765762
$ws.Range("M136").Value = -16231827
$ws.Range("N136").Value = -3770862

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 299
$ws.Range("I4").Value = 315.27274
$ws.Range("K4").Value = 315.27274
$ws.Range("M4").Value = -200.27274
$ws.Range("H21").Value = 34996.5
$ws.Range("J21").Value = 34996.5
$ws.Range("L21").Value = 34996.5
$ws.Range("N21").Value = -35468.5
$ws.Range("H86").Value = 1932.2727
$ws.Range("I86").Value = 1459.9375
$ws.Range("K86").Value = 1459.9375
$ws.Range("M86").Value = -336.9375
$ws.Range("H89").Value = 1932.2727
$ws.Range("I89").Value = 1459.9375
$ws.Range("K89").Value = 7299.6875
$ws.Range("M89").Value = -1683.6875
$ws.Range("H134").Value = 5265651
$ws.Range("I134").Value = 2585.9375
$ws.Range("K134").Value = 7757.8125
$ws.Range("M134").Value = -5222.8125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24052636
$ws.Range("I31").Value = 31252372
$ws.Range("K31").Value = 31252372
$ws.Range("M31").Value = -31252077
$ws.Range("H34").Value = 24052636
$ws.Range("I34").Value = 31252372
$ws.Range("K34").Value = 31252372
$ws.Range("M34").Value = -31252170
$ws.Range("H58").Value = 2839.5
$ws.Range("I58").Value = 2816.2917
$ws.Range("J58").Value = 2909.125
$ws.Range("K58").Value = 2816.2917
$ws.Range("L58").Value = 2909.125
$ws.Range("M58").Value = -2613.2917
$ws.Range("N58").Value = -3315.125
$ws.Range("H62").Value = 6907.8887
$ws.Range("J62").Value = 8010.3335
$ws.Range("L62").Value = 8010.3335
$ws.Range("N62").Value = -9258.333500000001
$ws.Range("H65").Value = 6907.8887
$ws.Range("J65").Value = 8010.3335
$ws.Range("L65").Value = 40051.6675
$ws.Range("N65").Value = -46291.6675
$ws.Range("H132").Value = 3294.4092
$ws.Range("I132").Value = 2799.125
$ws.Range("J132").Value = 4615.1665
$ws.Range("K132").Value = 8397.375
$ws.Range("L132").Value = 13845.4995
$ws.Range("M132").Value = -5867.375
$ws.Range("N132").Value = -18905.4995
$ws.Range("H134").Value = 3040.4583
$ws.Range("I134").Value = 3172.3157
$ws.Range("J134").Value = 2539.4
$ws.Range("K134").Value = 9516.947100000001
$ws.Range("L134").Value = 7618.200000000001
$ws.Range("M134").Value = -6981.947100000001
$ws.Range("N134").Value = -12688.2
$ws.Range("H136").Value = 2839.5
$ws.Range("I136").Value = 2816.2917
$ws.Range("J136").Value = 2909.125
$ws.Range("K136").Value = 8448.875100000001
$ws.Range("L136").Value = 8727.375
$ws.Range("M136").Value = -5898.875100000001
$ws.Range("N136").Value = -13827.375

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5555921
$ws.Range("I4").Value = 6250361
$ws.Range("K4").Value = 18751083
$ws.Range("M4").Value = -18750971
$ws.Range("H129").Value = 6340.846
$ws.Range("I129").Value = 3410.8
$ws.Range("J129").Value = 16107.667
$ws.Range("K129").Value = 10232.4
$ws.Range("L129").Value = 48323.001
$ws.Range("M129").Value = -5232.400000000001
$ws.Range("N129").Value = -58323.001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 508.14285
$ws.Range("I107").Value = 528.75
$ws.Range("J107").Value = 480.66666
$ws.Range("K107").Value = 528.75
$ws.Range("L107").Value = 480.66666
$ws.Range("M107").Value = 1391.25
$ws.Range("N107").Value = -4320.66666
$ws.Range("H132").Value = 4350069.5
$ws.Range("I132").Value = 2205.8333
$ws.Range("K132").Value = 6617.499899999999
$ws.Range("M132").Value = -4087.499899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 28998.25
$ws.Range("J23").Value = 28998.25
$ws.Range("L23").Value = 28998.25
$ws.Range("N23").Value = -29458.25
$ws.Range("H55").Value = 1827.2142
$ws.Range("I55").Value = 1574.2
$ws.Range("J55").Value = 1967.7778
$ws.Range("K55").Value = 1574.2
$ws.Range("L55").Value = 1967.7778
$ws.Range("M55").Value = -1401.2
$ws.Range("N55").Value = -2313.7778
$ws.Range("H68").Value = 4632173
$ws.Range("J68").Value = 3465
$ws.Range("L68").Value = 3465
$ws.Range("N68").Value = -4963
$ws.Range("H71").Value = 4632173
$ws.Range("J71").Value = 3465
$ws.Range("L71").Value = 17325
$ws.Range("N71").Value = -24813
$ws.Range("H94").Value = 75613.86
$ws.Range("J94").Value = 75613.86
$ws.Range("L94").Value = 75613.86
$ws.Range("N94").Value = -76965.86
$ws.Range("H132").Value = 4357.1665
$ws.Range("I132").Value = 2474.7
$ws.Range("J132").Value = 6710.25
$ws.Range("K132").Value = 7424.099999999999
$ws.Range("L132").Value = 20130.75
$ws.Range("M132").Value = -4894.099999999999
$ws.Range("N132").Value = -25190.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10192.917
$ws.Range("J62").Value = 12456.167
$ws.Range("L62").Value = 12456.167
$ws.Range("N62").Value = -13704.167
$ws.Range("H65").Value = 10192.917
$ws.Range("J65").Value = 12456.167
$ws.Range("L65").Value = 62280.835
$ws.Range("N65").Value = -68520.83499999999
$ws.Range("H136").Value = 181622.27
$ws.Range("I136").Value = 3185.8113
$ws.Range("J136").Value = 3333999.8
$ws.Range("K136").Value = 9557.4339
$ws.Range("L136").Value = 10001999.4
$ws.Range("M136").Value = -7007.4339
$ws.Range("N136").Value = -10007099.4
